# Auto-generated edit script: updates FY "leve profit" projection
# figures across several sheets following a scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 897.5526
$ws.Range("I98").Value = 507.21213
$ws.Range("J98").Value = 3473.8
$ws.Range("K98").Value = 507.21213
$ws.Range("L98").Value = 3473.8
$ws.Range("M98").Value = 990.78787
$ws.Range("N98").Value = -6469.8
$ws.Range("H122").Value = 897.5526
$ws.Range("I122").Value = 507.21213
$ws.Range("J122").Value = 3473.8
$ws.Range("K122").Value = 1521.63639
$ws.Range("L122").Value = 10421.4
$ws.Range("M122").Value = 928.3636099999999
$ws.Range("N122").Value = -15321.4
$ws.Range("H132").Value = 2610.9622
$ws.Range("I132").Value = 1619.8718
$ws.Range("J132").Value = 5371.857
$ws.Range("K132").Value = 4859.6154
$ws.Range("L132").Value = 16115.571
$ws.Range("M132").Value = -2329.6154
$ws.Range("N132").Value = -21175.571
$ws.Range("H138").Value = 3021.8064
$ws.Range("I138").Value = 2104.7693
$ws.Range("J138").Value = 3684.111
$ws.Range("K138").Value = 6314.3079
$ws.Range("L138").Value = 11052.333
$ws.Range("M138").Value = -1174.3079
$ws.Range("N138").Value = -21332.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1485.5883
$ws.Range("I97").Value = 1444.4546
$ws.Range("J97").Value = 1561
$ws.Range("K97").Value = 1444.4546
$ws.Range("L97").Value = 1561
$ws.Range("M97").Value = -948.4546
$ws.Range("N97").Value = -2553
$ws.Range("H113").Value = 49922.5
$ws.Range("J113").Value = 49922.5
$ws.Range("L113").Value = 49922.5
$ws.Range("N113").Value = -58600.5
$ws.Range("H122").Value = 1883.1305
$ws.Range("I122").Value = 1810.6
$ws.Range("J122").Value = 2366.6667
$ws.Range("K122").Value = 5431.799999999999
$ws.Range("L122").Value = 7100.000100000001
$ws.Range("M122").Value = -2981.799999999999
$ws.Range("N122").Value = -12000.0001
$ws.Range("H132").Value = 92533.62
$ws.Range("I132").Value = 144273.06
$ws.Range("J132").Value = 1989.6
$ws.Range("K132").Value = 432819.18
$ws.Range("L132").Value = 5968.799999999999
$ws.Range("M132").Value = -430289.18
$ws.Range("N132").Value = -11028.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1743.9474
$ws.Range("I86").Value = 1832.3334
$ws.Range("J86").Value = 1412.5
$ws.Range("K86").Value = 1832.3334
$ws.Range("L86").Value = 1412.5
$ws.Range("M86").Value = -709.3334
$ws.Range("N86").Value = -3658.5
$ws.Range("H89").Value = 1743.9474
$ws.Range("I89").Value = 1832.3334
$ws.Range("J89").Value = 1412.5
$ws.Range("K89").Value = 9161.666999999999
$ws.Range("L89").Value = 7062.5
$ws.Range("M89").Value = -3545.666999999999
$ws.Range("N89").Value = -18294.5
$ws.Range("H97").Value = 10276.8
$ws.Range("I97").Value = 7846
$ws.Range("K97").Value = 7846
$ws.Range("M97").Value = -6855
$ws.Range("H99").Value = 2771
$ws.Range("I99").Value = 1670
$ws.Range("J99").Value = 3242.8572
$ws.Range("K99").Value = 1670
$ws.Range("L99").Value = 3242.8572
$ws.Range("M99").Value = -172
$ws.Range("N99").Value = -6238.8572
$ws.Range("H134").Value = 76293
$ws.Range("I134").Value = 86966.46000000001
$ws.Range("K134").Value = 260899.38
$ws.Range("M134").Value = -258364.38

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1883.3112
$ws.Range("I31").Value = 1541.2894
$ws.Range("J31").Value = 3740
$ws.Range("K31").Value = 1541.2894
$ws.Range("L31").Value = 3740
$ws.Range("M31").Value = -1246.2894
$ws.Range("N31").Value = -4330
$ws.Range("H34").Value = 1883.3112
$ws.Range("I34").Value = 1541.2894
$ws.Range("J34").Value = 3740
$ws.Range("K34").Value = 1541.2894
$ws.Range("L34").Value = 3740
$ws.Range("M34").Value = -1339.2894
$ws.Range("N34").Value = -4144
$ws.Range("H86").Value = 58830640
$ws.Range("I86").Value = 83342290
$ws.Range("J86").Value = 2680
$ws.Range("K86").Value = 83342290
$ws.Range("L86").Value = 2680
$ws.Range("M86").Value = -83341167
$ws.Range("N86").Value = -4926
$ws.Range("H89").Value = 58830640
$ws.Range("I89").Value = 83342290
$ws.Range("J89").Value = 2680
$ws.Range("K89").Value = 416711450
$ws.Range("L89").Value = 13400
$ws.Range("M89").Value = -416705834
$ws.Range("N89").Value = -24632
$ws.Range("H98").Value = 52500
$ws.Range("J98").Value = 52500
$ws.Range("L98").Value = 52500
$ws.Range("N98").Value = -56992
$ws.Range("H99").Value = 1703.7368
$ws.Range("I99").Value = 1624.0834
$ws.Range("J99").Value = 1840.2858
$ws.Range("K99").Value = 1624.0834
$ws.Range("L99").Value = 1840.2858
$ws.Range("M99").Value = -126.0834
$ws.Range("N99").Value = -4836.2858
$ws.Range("H126").Value = 1703.7368
$ws.Range("I126").Value = 1624.0834
$ws.Range("J126").Value = 1840.2858
$ws.Range("K126").Value = 4872.2502
$ws.Range("L126").Value = 5520.857400000001
$ws.Range("M126").Value = -2402.2502
$ws.Range("N126").Value = -10460.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 710.9583
$ws.Range("I113").Value = 671.8889
$ws.Range("J113").Value = 734.4
$ws.Range("K113").Value = 2015.6667
$ws.Range("L113").Value = 2203.2
$ws.Range("M113").Value = 154.3332999999998
$ws.Range("N113").Value = -6543.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3530
$ws.Range("I97").Value = 3530
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3530
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3034
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5265119
$ws.Range("I7").Value = 7144439.5
$ws.Range("J7").Value = 3020.6
$ws.Range("K7").Value = 7144439.5
$ws.Range("L7").Value = 3020.6
$ws.Range("M7").Value = -7144327.5
$ws.Range("N7").Value = -3244.6
$ws.Range("H61").Value = 2392.0833
$ws.Range("I61").Value = 1800
$ws.Range("J61").Value = 3221
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 3221
$ws.Range("M61").Value = -1598
$ws.Range("N61").Value = -3625
$ws.Range("H113").Value = 2392.0833
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 3221
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 3221
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7561
$ws.Range("H122").Value = 3073.3794
$ws.Range("I122").Value = 2468.625
$ws.Range("J122").Value = 3817.6924
$ws.Range("K122").Value = 7405.875
$ws.Range("L122").Value = 11453.0772
$ws.Range("M122").Value = -4955.875
$ws.Range("N122").Value = -16353.0772
$ws.Range("H126").Value = 5265119
$ws.Range("I126").Value = 7144439.5
$ws.Range("J126").Value = 3020.6
$ws.Range("K126").Value = 21433318.5
$ws.Range("L126").Value = 9061.799999999999
$ws.Range("M126").Value = -21430848.5
$ws.Range("N126").Value = -14001.8
$ws.Range("H133").Value = 48554
$ws.Range("J133").Value = 48554
$ws.Range("L133").Value = 48554
$ws.Range("N133").Value = -53614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 48714.5
$ws.Range("J46").Value = 48714.5
$ws.Range("L46").Value = 48714.5
$ws.Range("N46").Value = -49176.5
$ws.Range("H134").Value = 48714.5
$ws.Range("J134").Value = 48714.5
$ws.Range("L134").Value = 146143.5
$ws.Range("N134").Value = -151213.5
